$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The only data change in the sheet is cell A2: "CasesTab" -> "ParticipantsTab"
$ws.Range("A2").Value = "ParticipantsTab"

# Update the selected cell/active cell from B3 to A2
$ws.Range("A2").Select()
